# "Add produce Logic of SLG Building"
# The Building struct definition sheet ("Property") had its explicit "ID" field
# row removed, and the type of a few fields (Prefab, NormalStateFunc,
# UpStateFunc, Desc) was fixed from "int" to "string".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the field type for the bottom four data rows (currently rows 5-8,
# i.e. Prefab/NormalStateFunc/UpStateFunc/Desc) from "int" to "string"
# before the ID row is removed and everything shifts up.
$ws.Range("B5:B8").Value = "string"

# Remove the "ID" field row (row 2); the remaining rows shift up so the
# struct now starts with "Type".
$ws.Rows(2).Delete()

# The list data-validation that applied below the data table needs to
# keep covering everything below the (now one-row-shorter) table.
$ws.Range("F8:F1048576").Validation.Delete()
$ws.Range("F8:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the sheet's recorded selection/active cell.
$ws.Range("G14").Select() | Out-Null
